$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TitleSheet")

$ws.Range("C1").Value = "ID"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2

$ws.Range("C3").Select()
